$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAccounts = @(
    @("Shaun_M70@gmail.com", 47144029),
    @("J.Biesheuvel56@yahoo.com", 20181690),
    @("Mohsin.S@yahoo.com", 26571302),
    @("J_Vente@hotmail.com", 40987142),
    @("Pia_A@yahoo.com", 58075440),
    @("Tooske_W@xs4all.nl", 68073468),
    @("Badia.G81@hotmail.com", 79961625),
    @("ChanineG@xs4all.nl", 89579070),
    @("Dila.Stein@kpnmail.nl", 17167209)
)

$startRow = 38
for ($i = 0; $i -lt $newAccounts.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newAccounts[$i][0]
    $ws.Cells.Item($row, 2).Value = $newAccounts[$i][1]
}
